$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name / title
$ws.Name = "Through 2021-11-08"

# Update label for November row
$ws.Range("A12").Value = "November (through 11-08)"

# Update November (row 12) values
$ws.Range("B12").Value = 10
$ws.Range("C12").Value = 19
$ws.Range("D12").Value = 29
$ws.Range("E12").Value = 20
$ws.Range("F12").Value = 12
$ws.Range("G12").Value = 50
$ws.Range("H12").Value = 55

# Update Total (row 13) values
$ws.Range("B13").Value = 268
$ws.Range("C13").Value = 505
$ws.Range("D13").Value = 739
$ws.Range("E13").Value = 635
$ws.Range("F13").Value = 494
$ws.Range("G13").Value = 1107
$ws.Range("H13").Value = 1499
